$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text/value assignments (names, links, percent changes, non-numeric-looking prices) ---
$ws.Range("D2").Value = "66.394.71"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.336.97"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "3.329.78"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "3.870.90"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.445.44"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "3.338.65"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +10.77%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +9.92%  "
$ws.Range("E33").Value = "  +4.14%  "
$ws.Range("E34").Value = "  +8.91%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.724.53"
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E41").Value = "  +6.00%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E42").Value = "  +7.60%  "
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("E43").Value = "  +12.39%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +0.17%  "

# --- Numeric-looking price strings: must stay text, so set via formula then paste-special as values ---
$ws.Range("D5").Formula = "=""189.21"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""562.52"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""0.592"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = "=""48.00"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("D14").Formula = "=""8.71"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D16").Formula = "=""608.16"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("D18").Formula = "=""18.13"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D21").Formula = "=""11.19"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D23").Formula = "=""18.65"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D24").Formula = "=""5.17"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=""100.78"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D26").Formula = "=""4.02"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""6.03"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""2.78"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""9.74"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D30").Formula = "=""8.71"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D31").Formula = "=""30.63"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""6.86"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""3.90"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D34").Formula = "=""581.97"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("D35").Formula = "=""11.18"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D36").Formula = "=""0.106"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D38").Formula = "=""57.31"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D39").Formula = "=""1.00"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D41").Formula = "=""0.132"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=""34.26"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D43").Formula = "=""3.48"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=""3.31"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D45").Formula = "=""2.72"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=""3.34"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D49").Formula = "=""0.131"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)

$excel.CutCopyMode = 0
